$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells (I1, J1) so they pick up style index 1 (bold, bordered,
# centered) just like the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-57: column I ("I0") and column J ("IF") values
$values = @(
    @(2, 7, 8),
    @(3, 8, 8),
    @(4, 8, 8),
    @(5, 8, 8),
    @(6, 6, 7),
    @(7, 8, 8),
    @(8, 8, 8),
    @(9, 7, 7),
    @(10, 7, 7),
    @(11, 8, 8),
    @(12, 6, 7),
    @(13, 7, 7),
    @(14, 8, 9),
    @(15, 8, 8),
    @(16, 7, 7),
    @(17, 6, 6),
    @(18, 6, 6),
    @(19, 8, 8),
    @(20, 7, 7),
    @(21, 7, 7),
    @(22, 8, 8),
    @(23, 8, 8),
    @(24, 4, 5),
    @(25, 8, 8),
    @(26, 6, 7),
    @(27, 7, 7),
    @(28, 9, 9),
    @(29, 8, 8),
    @(30, 6, 7),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 7, 7),
    @(34, 9, 9),
    @(35, 7, 7),
    @(36, 7, 8),
    @(37, 8, 8),
    @(38, 8, 8),
    @(39, 10, 10),
    @(40, 6, 6),
    @(41, 7, 7),
    @(42, 8, 9),
    @(43, 8, 8),
    @(44, 8, 9),
    @(45, 8, 8),
    @(46, 7, 8),
    @(47, 9, 9),
    @(48, 10, 10),
    @(49, 7, 8),
    @(50, 7, 7),
    @(51, 7, 7),
    @(52, 6, 7),
    @(53, 9, 9),
    @(54, 8, 8),
    @(55, 6, 6),
    @(56, 6, 6),
    @(57, 4, 4)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
